# Auto-generated edit script applying value changes per the commit diff.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC (12 cell updates) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H39").Value = 1056.3077
$ws.Range("I39").Value = 203
$ws.Range("K39").Value = 609
$ws.Range("M39").Value = -313
$ws.Range("H70").Value = 1357047.9
$ws.Range("J70").Value = 2635
$ws.Range("L70").Value = 7905
$ws.Range("N70").Value = -8445
$ws.Range("H73").Value = 1357047.9
$ws.Range("J73").Value = 2635
$ws.Range("L73").Value = 7905
$ws.Range("N73").Value = -9777

# --- Sheet: ARM (32 cell updates) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 8404494
$ws.Range("I61").Value = 9528717
$ws.Range("J61").Value = 2502325
$ws.Range("K61").Value = 9528717
$ws.Range("L61").Value = 2502325
$ws.Range("M61").Value = -9528505
$ws.Range("N61").Value = -2502749
$ws.Range("H88").Value = 2028.7
$ws.Range("I88").Value = 2550
$ws.Range("J88").Value = 1970.7778
$ws.Range("K88").Value = 2550
$ws.Range("L88").Value = 1970.7778
$ws.Range("M88").Value = -2144
$ws.Range("N88").Value = -2782.7778
$ws.Range("H91").Value = 2028.7
$ws.Range("I91").Value = 2550
$ws.Range("J91").Value = 1970.7778
$ws.Range("K91").Value = 2550
$ws.Range("L91").Value = 1970.7778
$ws.Range("M91").Value = -1146
$ws.Range("N91").Value = -4778.7778
$ws.Range("H132").Value = 3126996
$ws.Range("I132").Value = 1877.9166
$ws.Range("K132").Value = 5633.7498
$ws.Range("M132").Value = -3103.7498
$ws.Range("H136").Value = 8404494
$ws.Range("I136").Value = 9528717
$ws.Range("J136").Value = 2502325
$ws.Range("K136").Value = 28586151
$ws.Range("L136").Value = 7506975
$ws.Range("M136").Value = -28583601
$ws.Range("N136").Value = -7512075

# --- Sheet: BSM (46 cell updates) ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H21").Value = 34999.5
$ws.Range("J21").Value = 34999.5
$ws.Range("L21").Value = 34999.5
$ws.Range("N21").Value = -35471.5
$ws.Range("H26").Value = 25297.666
$ws.Range("I26").Value = 10357.4
$ws.Range("J26").Value = 99999
$ws.Range("K26").Value = 10357.4
$ws.Range("L26").Value = 99999
$ws.Range("M26").Value = -10065.4
$ws.Range("N26").Value = -100583
$ws.Range("H76").Value = 46166.168
$ws.Range("J76").Value = 46166.168
$ws.Range("L76").Value = 46166.168
$ws.Range("N76").Value = -46796.168
$ws.Range("H79").Value = 46166.168
$ws.Range("J79").Value = 46166.168
$ws.Range("L79").Value = 46166.168
$ws.Range("N79").Value = -48350.168
$ws.Range("H86").Value = 71240.10000000001
$ws.Range("I86").Value = 87000.125
$ws.Range("K86").Value = 87000.125
$ws.Range("M86").Value = -85877.125
$ws.Range("H89").Value = 71240.10000000001
$ws.Range("I89").Value = 87000.125
$ws.Range("K89").Value = 435000.625
$ws.Range("M89").Value = -429384.625
$ws.Range("H100").Value = 30838.4
$ws.Range("J100").Value = 30838.4
$ws.Range("L100").Value = 30838.4
$ws.Range("N100").Value = -33002.4
$ws.Range("H101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("N101").Value = ""
$ws.Range("H102").Value = 33199.75
$ws.Range("I102").Value = 10933.333
$ws.Range("K102").Value = 10933.333
$ws.Range("M102").Value = -7688.333000000001
$ws.Range("H107").Value = 4840
$ws.Range("I107").Value = 5780
$ws.Range("J107").Value = 1550
$ws.Range("K107").Value = 5780
$ws.Range("L107").Value = 1550
$ws.Range("M107").Value = -3860
$ws.Range("N107").Value = -5390

# --- Sheet: CRP (52 cell updates) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 766.6667
$ws.Range("I22").Value = 650.25
$ws.Range("J22").Value = 999.5
$ws.Range("K22").Value = 650.25
$ws.Range("L22").Value = 999.5
$ws.Range("M22").Value = -300.25
$ws.Range("N22").Value = -1699.5
$ws.Range("H23").Value = 34990
$ws.Range("I23").Value = 34990
$ws.Range("K23").Value = 34990
$ws.Range("M23").Value = -34750
$ws.Range("H27").Value = 34990
$ws.Range("I27").Value = 34990
$ws.Range("K27").Value = 34990
$ws.Range("M27").Value = -34798
$ws.Range("H31").Value = 31275078
$ws.Range("I31").Value = 37065892
$ws.Range("J31").Value = 4682.8
$ws.Range("K31").Value = 37065892
$ws.Range("L31").Value = 4682.8
$ws.Range("M31").Value = -37065597
$ws.Range("N31").Value = -5272.8
$ws.Range("H34").Value = 31275078
$ws.Range("I34").Value = 37065892
$ws.Range("J34").Value = 4682.8
$ws.Range("K34").Value = 37065892
$ws.Range("L34").Value = 4682.8
$ws.Range("M34").Value = -37065690
$ws.Range("N34").Value = -5086.8
$ws.Range("H58").Value = 3317.4
$ws.Range("I58").Value = 3123.1304
$ws.Range("K58").Value = 3123.1304
$ws.Range("M58").Value = -2920.1304
$ws.Range("H94").Value = 2465
$ws.Range("I94").Value = 1780.5
$ws.Range("K94").Value = 1780.5
$ws.Range("M94").Value = -1329.5
$ws.Range("H122").Value = 2426.25
$ws.Range("I122").Value = 2317.1667
$ws.Range("J122").Value = 2753.5
$ws.Range("K122").Value = 6951.500100000001
$ws.Range("L122").Value = 8260.5
$ws.Range("M122").Value = -4501.500100000001
$ws.Range("N122").Value = -13160.5
$ws.Range("H132").Value = 2441.8696
$ws.Range("I132").Value = 2091.5625
$ws.Range("K132").Value = 6274.6875
$ws.Range("M132").Value = -3744.6875
$ws.Range("H136").Value = 3317.4
$ws.Range("I136").Value = 3123.1304
$ws.Range("K136").Value = 9369.3912
$ws.Range("M136").Value = -6819.3912

# --- Sheet: CUL (32 cell updates) ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 524.7895
$ws.Range("I5").Value = 187.93333
$ws.Range("K5").Value = 563.79999
$ws.Range("M5").Value = -451.79999
$ws.Range("H86").Value = 499.5
$ws.Range("J86").Value = 499
$ws.Range("L86").Value = 1497
$ws.Range("N86").Value = -3869
$ws.Range("H89").Value = 499.5
$ws.Range("J89").Value = 499
$ws.Range("L89").Value = 4491
$ws.Range("N89").Value = -16347
$ws.Range("H108").Value = 16680
$ws.Range("I108").Value = 27
$ws.Range("K108").Value = 81
$ws.Range("M108").Value = 2799
$ws.Range("H112").Value = 17958.25
$ws.Range("I112").Value = 8500
$ws.Range("K112").Value = 25500
$ws.Range("M112").Value = -24392
$ws.Range("H120").Value = 30807.25
$ws.Range("I120").Value = 29965.334
$ws.Range("K120").Value = 89896.00199999999
$ws.Range("M120").Value = -85058.00199999999
$ws.Range("H133").Value = 45965.668
$ws.Range("I133").Value = 39664.145
$ws.Range("K133").Value = 118992.435
$ws.Range("M133").Value = -113932.435
$ws.Range("H135").Value = 524.7895
$ws.Range("I135").Value = 187.93333
$ws.Range("K135").Value = 1691.39997
$ws.Range("M135").Value = 843.6000299999998

# --- Sheet: GSM (11 cell updates) ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 16303950
$ws.Range("I126").Value = 23289638
$ws.Range("K126").Value = 69868914
$ws.Range("M126").Value = -69866444
$ws.Range("H132").Value = 3746740.8
$ws.Range("I132").Value = 3377.805
$ws.Range("J132").Value = 19094528
$ws.Range("K132").Value = 10133.415
$ws.Range("L132").Value = 57283584
$ws.Range("M132").Value = -7603.414999999999
$ws.Range("N132").Value = -57288644

# --- Sheet: LTW (12 cell updates) ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 5212358
$ws.Range("J68").Value = 6751.5
$ws.Range("L68").Value = 6751.5
$ws.Range("N68").Value = -8249.5
$ws.Range("H71").Value = 5212358
$ws.Range("J71").Value = 6751.5
$ws.Range("L71").Value = 33757.5
$ws.Range("N71").Value = -41245.5
$ws.Range("H132").Value = 3901.0232
$ws.Range("I132").Value = 2540.7273
$ws.Range("K132").Value = 7622.1819
$ws.Range("M132").Value = -5092.1819

# --- Sheet: WVR (8 cell updates) ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 24450
$ws.Range("H132").Value = 225049.6
$ws.Range("I132").Value = 2706.5405
$ws.Range("J132").Value = 1253386.2
$ws.Range("K132").Value = 8119.6215
$ws.Range("L132").Value = 3760158.6
$ws.Range("M132").Value = -5589.6215
$ws.Range("N132").Value = -3765218.6

